$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 3, 5, 6, 7, 9, 10, 11, 12, 14)

$data = @(
    @(1.300208399224204, 0.04395605117662171, 0.08576357180315597, 3.109416300997921, 0.002580363329603541, 1.897989459481202, 0.1473016881408207, 1.052472831731592, 0.3802095517735609, 3.087961604294925),
    @(1.262008088567597, 0.03900796721909217, 0.08536323208581642, 3.099647560060859, 0.00258421229339725, 1.90088666723711, 0.148000278681635, 1.012499874882309, 0.3740354005221178, 3.107356460029003),
    @(1.239253311982111, 0.03598554884227667, 0.08514649481235281, 3.095007892897868, 0.002586701312871189, 1.903400244391122, 0.1484465133231436, 0.9885536846202285, 0.3704258763042532, 3.120062264798023),
    @(1.23015696916417, 0.03475774506732421, 0.08506551445782407, 3.093458891801617, 0.00258774732612716, 1.904609283869462, 0.1486327128030769, 0.9789457033002122, 0.3690007032537608, 3.125440397108541),
    @(1.228657192275563, 0.03455409923994068, 0.08505251208134901, 3.093222323065035, 0.002587922934600121, 1.904821201502529, 0.1486638943003269, 0.977359384510379, 0.3687668201958019, 3.126345534384463),
    @(1.23912992083234, 0.0359689747674139, 0.08514537291426549, 3.09498561883305, 0.002586715291149346, 1.903415801942494, 0.1484490068329833, 0.9884234994248686, 0.3704064705783168, 3.12013398480866),
    @(1.286891733980127, 0.04224662756631403, 0.08561951898880515, 3.10576617142803, 0.002581664414898588, 1.898835927220844, 0.1475389756223828, 1.038566199949713, 0.3780431149039742, 3.094483432505527),
    @(1.386103088939024, 0.05468688258913801, 0.08677874643219496, 3.137684189206269, 0.002572752814401941, 1.895685977374804, 0.1458913633058292, 1.141639203815913, 0.3944542640890916, 3.050511739753432),
    @(1.462378540704947, 0.06391369311151607, 0.08776878907611874, 3.167710623018678, 0.002566804549532065, 1.896931374598601, 0.1447639082748622, 1.220273187720721, 0.4073834165413075, 3.022067590993757),
    @(1.497814147025906, 0.06813183112384991, 0.08824892185359801, 3.182800251075591, 0.002564227261150116, 1.898272134787334, 0.1442689372465011, 1.256680693373909, 0.4134539561867854, 3.009967015442712),
    @(1.511338616458772, 0.06973224555716229, 0.08843498601315503, 3.188720023782636, 0.002563269699776612, 1.898891248222668, 0.1440840748175591, 1.270558906367626, 0.4157798040180722, 3.005505567155495),
    @(1.508421180050846, 0.06938742826963562, 0.088394725311705, 3.187435948691331, 0.002563475110889893, 1.898752955613119, 0.144123773899083, 1.267565916987252, 0.4152776890757792, 3.006461045725388),
    @(1.498924694978598, 0.06826343551544767, 0.08826414449758246, 3.183283152855211, 0.002564148113728116, 1.898320836688896, 0.1442536769694116, 1.257820629163348, 0.4136447632808711, 3.009597547940132),
    @(1.493121590373846, 0.06757536392579766, 0.08818471232124381, 3.180766227152304, 0.002564562741148414, 1.898070660461343, 0.1443335812426598, 1.251863271098586, 0.412648070583657, 3.011534479585706),
    @(1.460077557240083, 0.06363845559958747, 0.08773800743013282, 3.166753261242235, 0.002566975561634291, 1.896859337543418, 0.1447966159543479, 1.217906675573971, 0.4069904861188434, 3.022875284406155),
    @(1.439994809148288, 0.06122869495291638, 0.08747156808256662, 3.158523104450538, 0.002568488622604173, 1.896314572078609, 0.1450852582938658, 1.197238393913494, 0.4035680694983768, 3.030047450751425),
    @(1.428513179015965, 0.05984461738842128, 0.08732112195728803, 3.15392397896872, 0.002569371005851682, 1.896074106844985, 0.1452529640811067, 1.185410459129997, 0.4016173825355622, 3.034251634221619),
    @(1.424637630548432, 0.05937632478512, 0.08727066579758258, 3.152389920455136, 0.002569671848271271, 1.896005202943378, 0.1453100361253936, 1.181416015522871, 0.4009599734537517, 3.035688652357308),
    @(1.442125467705694, 0.06148501520682714, 0.08749964115671105, 3.15938528393329, 0.002568326301950754, 1.896365021112651, 0.1450543573195868, 1.19943236725544, 0.4039305505297932, 3.029275789948287),
    @(1.501711174366903, 0.06859349454721553, 0.08830238419508873, 3.184497348448048, 0.002563949937846144, 1.89844473669072, 0.1442154515041119, 1.260680574033557, 0.4141236594912243, 3.008673002080172),
    @(1.541270285428254, 0.07325739767478012, 0.08885177132659194, 3.202108251715146, 0.00256119694913944, 1.90045328792155, 0.1436821701527755, 1.301243018019761, 0.4209431564984953, 2.99591194635174),
    @(1.520100533443326, 0.07076649409825109, 0.08855629874712889, 3.192599310328561, 0.002562656490443996, 1.899321849950496, 0.1439654218303019, 1.279545305453297, 0.4172890689200983, 3.002658287648828),
    @(1.441161996942412, 0.0613691286993685, 0.08748694080397001, 3.158995079839855, 0.002568399648086331, 1.896341986553772, 0.1450683221616198, 1.198440302029923, 0.4037666201892733, 3.029624406066915),
    @(1.358669465864381, 0.05130676125932609, 0.08644073923036188, 3.127895114384572, 0.002575057974719138, 1.89591336729999, 0.1463224707607416, 1.113245928025833, 0.3898612614374031, 3.061729338462079)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($r, $cols[$j]).Value = $rowVals[$j]
    }
}

Write-Host "Done updating cells"